# Apply the "add testing for complex numbers" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing Duration values (column E) for rows 2-23 ---
$durations = @{
    2  = 0.0449395
    3  = 0.0027536
    4  = 0.0022373
    5  = 0.0026535
    6  = 0.0022807
    7  = 0.0078134
    8  = 0.0045891
    9  = 0.0044361
    10 = 0.0037888
    11 = 0.0060624
    12 = 0.0042094
    13 = 0.0049696
    14 = 0.0056153
    15 = 0.0044784
    16 = 0.0047477
    17 = 0.0037508
    18 = 0.004733
    19 = 0.0046678
    20 = 0.0048193
    21 = 0.0056896
    22 = 0.0051992
    23 = 0.0066391
}

foreach ($row in $durations.Keys) {
    $ws.Cells.Item($row, 5).Value = $durations[$row]
}

# --- Add new rows for the complex numbers tests ---
$newRows = @(
    @{ Row = 24; Name = "FactorielClassTest/testComplexNumbers(complexNumbers=value1)"; Duration = 0.0046625 },
    @{ Row = 25; Name = "FactorielClassTest/testComplexNumbers(complexNumbers=value2)"; Duration = 0.0040134 },
    @{ Row = 26; Name = "FactorielClassTest/testComplexNumbers(complexNumbers=value3)"; Duration = 0.0359294 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $true
    $ws.Cells.Item($r, 3).Value = $false
    $ws.Cells.Item($r, 4).Value = $false
    $ws.Cells.Item($r, 5).Value = $entry.Duration
}
